# Login Tests with Select
$wb = $excel.ActiveWorkbook

# Rename the second sheet ("SecondTest") to "addCustomerTest"
$ws2 = $wb.Worksheets.Item("SecondTest")
$ws2.Name = "addCustomerTest"

# Add a third column (postcode) of data to addCustomerTest and a 4th data row
$ws2.Range("C1").Value = "postcode"
$ws2.Range("C2").Value = 12345
$ws2.Range("C3").Value = 56789

$ws2.Range("A4").Value = "Tim"
$ws2.Range("B4").Value = "Fisher"
$ws2.Range("C4").Value = 45678

# Set selection on addCustomerTest (not the active tab afterwards)
$ws2.Activate()
$ws2.Range("G8").Select()

# Add a new sheet "openAccountTest" right after addCustomerTest (becomes the 3rd/last sheet)
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "openAccountTest"

$ws3.Range("A1").Value = "customer"
$ws3.Range("B1").Value = "currency"
$ws3.Range("A2").Value = "Joe Smith"
$ws3.Range("B2").Value = "Dollar"

# Make openAccountTest the active sheet / tab, with its own selection
$ws3.Activate()
$ws3.Range("B3").Select()
